$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) / Volume(1h) (E) columns, and the swapped Coin/Link (B/C)
# for rows 9 and 10. D/E cells are forced to Text format (NumberFormat "@")
# individually before assignment so values such as "1.00" / "19.30" keep their
# exact original text representation instead of being coerced to numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '72.439.83'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.646.81'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.85%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.02'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.55'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.36%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('B9').Value = 'LidoStakedEther'
$ws.Range('C9').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.646.66'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.89%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.172'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.97%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.130.31'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.92%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '72.288.51'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.92'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.672.55'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.13'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.97'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '371.55'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.15'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.24%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '71.01'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.25'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.65'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.782.75'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.03'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '496.55'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.99%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.45%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.51'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.30'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.113'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.12%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.25%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.48%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -6.02%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.327'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.38%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '152.87'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.61%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0749'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.88%  '

Write-Host "Applied cryptos list update"
